$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44901
$ws.Range("M2").Value = 150
$ws.Range("Q2").Value = '$/caja 10 kilos'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 1500
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44901
$ws.Range("K3").Value = 'Castle Brite'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 13000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 13000
$ws.Range("Q3").Value = '$/caja 10 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 1300
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44917
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 17000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("S4").Value = 944
$ws.Range("T4").Value = 18

# Row 5
$ws.Range("D5").Value = 44160
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("Q5").Value = '$/caja 15 kilos'
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("S5").Value = 1333
$ws.Range("T5").Value = 15

# Row 6
$ws.Range("D6").Value = 44189
$ws.Range("M6").Value = 50
$ws.Range("R6").Value = 'Provincia de San Felipe de Aconcagua'

# Row 7
$ws.Range("D7").Value = 44883
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 1400

# Row 8
$ws.Range("D8").Value = 44937
$ws.Range("K8").Value = 'Modesto'
$ws.Range("M8").Value = 230
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 15000
$ws.Range("Q8").Value = '$/caja 15 kilos'
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44910
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 240
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("Q9").Value = '$/caja 10 kilos'
$ws.Range("S9").Value = 1000
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44172
$ws.Range("L10").Value = 'Especial'
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = '$/caja 10 kilos'
$ws.Range("R10").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S10").Value = 1500
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 44172
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 250
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 11000
$ws.Range("P11").Value = 11000
$ws.Range("Q11").Value = '$/caja 10 kilos'
$ws.Range("R11").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S11").Value = 1100
$ws.Range("T11").Value = 10

# Row 12
$ws.Range("D12").Value = 44915
$ws.Range("K12").Value = 'Dina'
$ws.Range("M12").Value = 270
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44547
$ws.Range("M13").Value = 120
$ws.Range("Q13").Value = '$/caja 15 kilos'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1133
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = 44547
$ws.Range("M14").Value = 170
$ws.Range("Q14").Value = '$/caja 15 kilos'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 1000
$ws.Range("T14").Value = 15

# Row 15
$ws.Range("D15").Value = 44168
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 250
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("S15").Value = 1000

# Row 16
$ws.Range("D16").Value = 44168
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 17000
$ws.Range("O16").Value = 17000
$ws.Range("P16").Value = 17000
$ws.Range("Q16").Value = '$/caja 18 kilos'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 944
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44553
$ws.Range("L17").Value = 'Especial'
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 13000
$ws.Range("O17").Value = 13000
$ws.Range("P17").Value = 13000
$ws.Range("S17").Value = 1300

# Row 18
$ws.Range("D18").Value = 44553
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 120
$ws.Range("N18").Value = 11000
$ws.Range("O18").Value = 11000
$ws.Range("P18").Value = 11000
$ws.Range("S18").Value = 1100

# Row 19
$ws.Range("D19").Value = 44553
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 14000
$ws.Range("O19").Value = 14000
$ws.Range("P19").Value = 14000
$ws.Range("Q19").Value = '$/caja 15 kilos'
$ws.Range("S19").Value = 933
$ws.Range("T19").Value = 15

# Row 20
$ws.Range("D20").Value = 44900
$ws.Range("K20").Value = 'Castle Brite'
$ws.Range("N20").Value = 15000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 15545
$ws.Range("Q20").Value = '$/caja 10 kilos'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 1554
$ws.Range("T20").Value = 10

# Row 21
$ws.Range("D21").Value = 44186
$ws.Range("M21").Value = 150
$ws.Range("R21").Value = 'Región Metropolitana'

# Row 22
$ws.Range("D22").Value = 44566
$ws.Range("K22").Value = 'Modesto'
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 10000
$ws.Range("Q22").Value = '$/caja 10 kilos'
$ws.Range("S22").Value = 1000
$ws.Range("T22").Value = 10

# Row 23
$ws.Range("D23").Value = 44904
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 120
$ws.Range("Q23").Value = '$/bandeja 10 kilos'
$ws.Range("S23").Value = 1500
$ws.Range("T23").Value = 10

# Row 24
$ws.Range("D24").Value = 44904
$ws.Range("K24").Value = 'Castle Brite'
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 13000
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 13000
$ws.Range("S24").Value = 1300

# Row 25
$ws.Range("D25").Value = 44162
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 17000
$ws.Range("Q25").Value = '$/caja 16 kilos granel'
$ws.Range("R25").Value = 'Provincia de Limarí'
$ws.Range("S25").Value = 1062
$ws.Range("T25").Value = 16

# Row 26
$ws.Range("D26").Value = 44162
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 100
$ws.Range("Q26").Value = '$/caja 16 kilos granel'
$ws.Range("R26").Value = 'Provincia de Limarí'
$ws.Range("S26").Value = 938
$ws.Range("T26").Value = 16

# Row 27
$ws.Range("D27").Value = 44931
$ws.Range("K27").Value = 'Dina'
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 180
$ws.Range("N27").Value = 15000
$ws.Range("O27").Value = 15000
$ws.Range("P27").Value = 15000
$ws.Range("Q27").Value = '$/caja 16 kilos'
$ws.Range("S27").Value = 938
$ws.Range("T27").Value = 16

# Row 28
$ws.Range("D28").Value = 44902
$ws.Range("M28").Value = 560
$ws.Range("N28").Value = 14000
$ws.Range("O28").Value = 15000
$ws.Range("P28").Value = 14643
$ws.Range("R28").Value = 'Región de O''Higgins'
$ws.Range("S28").Value = 1464

# Row 29
$ws.Range("D29").Value = 44902
$ws.Range("L29").Value = 'Segunda'
$ws.Range("M29").Value = 190
$ws.Range("N29").Value = 13000
$ws.Range("O29").Value = 13000
$ws.Range("P29").Value = 13000
$ws.Range("Q29").Value = '$/caja 10 kilos'
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 1300
$ws.Range("T29").Value = 10

# Row 30
$ws.Range("D30").Value = 44161
$ws.Range("M30").Value = 150
$ws.Range("N30").Value = 20000
$ws.Range("O30").Value = 20000
$ws.Range("P30").Value = 20000
$ws.Range("Q30").Value = '$/caja 18 kilos granel'
$ws.Range("R30").Value = 'Provincia de Limarí'
$ws.Range("S30").Value = 1111
$ws.Range("T30").Value = 18

# Row 31
$ws.Range("D31").Value = 44181
$ws.Range("K31").Value = 'Dina'
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 220
$ws.Range("N31").Value = 17000
$ws.Range("O31").Value = 17000
$ws.Range("P31").Value = 17000
$ws.Range("Q31").Value = '$/caja 18 kilos'
$ws.Range("R31").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S31").Value = 944
$ws.Range("T31").Value = 18

# Row 32
$ws.Range("D32").Value = 44529
$ws.Range("N32").Value = 20000
$ws.Range("O32").Value = 20000
$ws.Range("P32").Value = 20000
$ws.Range("S32").Value = 1333

# Row 33
$ws.Range("D33").Value = 44187
$ws.Range("K33").Value = 'Dina'
$ws.Range("L33").Value = 'Primera'
$ws.Range("M33").Value = 120
$ws.Range("N33").Value = 16000
$ws.Range("O33").Value = 16000
$ws.Range("P33").Value = 16000
$ws.Range("Q33").Value = '$/caja 18 kilos'
$ws.Range("R33").Value = 'Provincia de Limarí'
$ws.Range("S33").Value = 889
$ws.Range("T33").Value = 18

# Row 34
$ws.Range("D34").Value = 44550
$ws.Range("L34").Value = 'Segunda'
$ws.Range("M34").Value = 200
$ws.Range("N34").Value = 15000
$ws.Range("O34").Value = 15000
$ws.Range("P34").Value = 15000
$ws.Range("S34").Value = 1000

# Row 35
$ws.Range("D35").Value = 44176
$ws.Range("L35").Value = 'Primera'
$ws.Range("M35").Value = 100
$ws.Range("N35").Value = 17000
$ws.Range("O35").Value = 17000
$ws.Range("P35").Value = 17000
$ws.Range("Q35").Value = '$/caja 18 kilos granel'
$ws.Range("R35").Value = 'Provincia de Limarí'
$ws.Range("S35").Value = 944
$ws.Range("T35").Value = 18

# Row 36
$ws.Range("D36").Value = 44543
$ws.Range("M36").Value = 100
$ws.Range("N36").Value = 18000
$ws.Range("O36").Value = 18000
$ws.Range("P36").Value = 18000
$ws.Range("Q36").Value = '$/caja 15 kilos'
$ws.Range("R36").Value = 'Región de O''Higgins'
$ws.Range("S36").Value = 1200
$ws.Range("T36").Value = 15

# Row 37
$ws.Range("D37").Value = 44543
$ws.Range("K37").Value = 'Castle Brite'
$ws.Range("L37").Value = 'Segunda'
$ws.Range("M37").Value = 50
$ws.Range("N37").Value = 15000
$ws.Range("O37").Value = 15000
$ws.Range("P37").Value = 15000
$ws.Range("Q37").Value = '$/caja 15 kilos'
$ws.Range("T37").Value = 15

# Row 38
$ws.Range("D38").Value = 44540
$ws.Range("K38").Value = 'Castle Brite'
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 120
$ws.Range("N38").Value = 20000
$ws.Range("O38").Value = 20000
$ws.Range("P38").Value = 20000
$ws.Range("Q38").Value = '$/caja 16 kilos'
$ws.Range("S38").Value = 1250
$ws.Range("T38").Value = 16

# Row 39
$ws.Range("L39").Value = 'Segunda'
$ws.Range("M39").Value = 200
$ws.Range("N39").Value = 15000
$ws.Range("O39").Value = 15000
$ws.Range("P39").Value = 15000
$ws.Range("S39").Value = 938

# Row 40
$ws.Range("D40").Value = 44907
$ws.Range("L40").Value = 'Primera'
$ws.Range("M40").Value = 120
$ws.Range("N40").Value = 20000
$ws.Range("O40").Value = 20000
$ws.Range("P40").Value = 20000
$ws.Range("Q40").Value = '$/caja 18 kilos'
$ws.Range("S40").Value = 1111
$ws.Range("T40").Value = 18

# Row 41
$ws.Range("D41").Value = 44533
$ws.Range("M41").Value = 170
$ws.Range("N41").Value = 18000
$ws.Range("O41").Value = 18000
$ws.Range("P41").Value = 18000
$ws.Range("Q41").Value = '$/caja 15 kilos'
$ws.Range("R41").Value = 'Región de O''Higgins'
$ws.Range("S41").Value = 1200
$ws.Range("T41").Value = 15

# Row 42
$ws.Range("D42").Value = 44533
$ws.Range("L42").Value = 'Segunda'
$ws.Range("M42").Value = 100
$ws.Range("N42").Value = 14000
$ws.Range("O42").Value = 14000
$ws.Range("P42").Value = 14000
$ws.Range("S42").Value = 933

# Row 43
$ws.Range("D43").Value = 44545
$ws.Range("L43").Value = 'Primera'
$ws.Range("M43").Value = 120
$ws.Range("N43").Value = 17000
$ws.Range("O43").Value = 17000
$ws.Range("P43").Value = 17000
$ws.Range("S43").Value = 1133

# Row 44
$ws.Range("D44").Value = 44545
$ws.Range("K44").Value = 'Castle Brite'
$ws.Range("L44").Value = 'Segunda'
$ws.Range("M44").Value = 100
$ws.Range("Q44").Value = '$/caja 15 kilos'
$ws.Range("S44").Value = 1000
$ws.Range("T44").Value = 15

# Row 45
$ws.Range("D45").Value = 44932
$ws.Range("K45").Value = 'Dina'
$ws.Range("M45").Value = 150
$ws.Range("N45").Value = 10000
$ws.Range("O45").Value = 10000
$ws.Range("P45").Value = 10000
$ws.Range("Q45").Value = '$/bandeja 10 kilos'
$ws.Range("T45").Value = 10

# Row 46
$ws.Range("D46").Value = 44567
$ws.Range("K46").Value = 'Modesto'
$ws.Range("M46").Value = 200
$ws.Range("N46").Value = 18000
$ws.Range("O46").Value = 18000
$ws.Range("P46").Value = 18000
$ws.Range("R46").Value = 'Región de O''Higgins'
$ws.Range("S46").Value = 1000

# Row 47
$ws.Range("D47").Value = 44567
$ws.Range("K47").Value = 'Modesto'
$ws.Range("L47").Value = 'Segunda'
$ws.Range("M47").Value = 250
$ws.Range("N47").Value = 13000
$ws.Range("O47").Value = 13000
$ws.Range("P47").Value = 13000
$ws.Range("S47").Value = 867

# Row 48
$ws.Range("D48").Value = 44179
$ws.Range("M48").Value = 150
$ws.Range("N48").Value = 18000
$ws.Range("O48").Value = 18000
$ws.Range("P48").Value = 18000
$ws.Range("R48").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S48").Value = 1000

# Row 49
$ws.Range("D49").Value = 44874
$ws.Range("K49").Value = 'Castle Brite'
$ws.Range("M49").Value = 60
$ws.Range("N49").Value = 30000
$ws.Range("O49").Value = 30000
$ws.Range("P49").Value = 30000
$ws.Range("Q49").Value = '$/bandeja 10 kilos'
$ws.Range("R49").Value = 'Provincia de Limarí'
$ws.Range("S49").Value = 3000
$ws.Range("T49").Value = 10

# Row 50
$ws.Range("D50").Value = 44918
$ws.Range("N50").Value = 18000
$ws.Range("O50").Value = 18000
$ws.Range("P50").Value = 18000
$ws.Range("Q50").Value = '$/caja 18 kilos'
$ws.Range("R50").Value = 'Región de O''Higgins'
$ws.Range("S50").Value = 1000
$ws.Range("T50").Value = 18

# Row 51
$ws.Range("D51").Value = 44530
$ws.Range("K51").Value = 'Castle Brite'
$ws.Range("M51").Value = 130
$ws.Range("N51").Value = 20000
$ws.Range("O51").Value = 20000
$ws.Range("P51").Value = 20000
$ws.Range("Q51").Value = '$/caja 15 kilos'
$ws.Range("S51").Value = 1333
$ws.Range("T51").Value = 15

# Row 52
$ws.Range("D52").Value = 44530
$ws.Range("L52").Value = 'Segunda'
$ws.Range("M52").Value = 150
$ws.Range("N52").Value = 15000
$ws.Range("O52").Value = 15000
$ws.Range("P52").Value = 15000
$ws.Range("Q52").Value = '$/caja 15 kilos'
$ws.Range("R52").Value = 'Región de O''Higgins'
$ws.Range("S52").Value = 1000
$ws.Range("T52").Value = 15

# Row 53
$ws.Range("D53").Value = 44911
$ws.Range("K53").Value = 'Dina'
$ws.Range("M53").Value = 250
$ws.Range("N53").Value = 20000
$ws.Range("O53").Value = 20000
$ws.Range("P53").Value = 20000
$ws.Range("Q53").Value = '$/caja 16 kilos'
$ws.Range("R53").Value = 'Región de O''Higgins'
$ws.Range("S53").Value = 1250
$ws.Range("T53").Value = 16

# Row 54
$ws.Range("D54").Value = 44911
$ws.Range("K54").Value = 'Dina'
$ws.Range("L54").Value = 'Primera'
$ws.Range("N54").Value = 18000
$ws.Range("O54").Value = 18000
$ws.Range("P54").Value = 18000
$ws.Range("Q54").Value = '$/caja 16 kilos'
$ws.Range("S54").Value = 1125
$ws.Range("T54").Value = 16

# Row 55
$ws.Range("D55").Value = 44174
$ws.Range("L55").Value = 'Especial'
$ws.Range("M55").Value = 200
$ws.Range("N55").Value = 15000
$ws.Range("O55").Value = 15000
$ws.Range("P55").Value = 15000
$ws.Range("R55").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S55").Value = 1500

# Row 56
$ws.Range("D56").Value = 44167
$ws.Range("L56").Value = 'Primera'
$ws.Range("M56").Value = 300
$ws.Range("N56").Value = 15000
$ws.Range("O56").Value = 15000
$ws.Range("P56").Value = 15000
$ws.Range("Q56").Value = '$/caja 16 kilos granel'
$ws.Range("R56").Value = 'Provincia de Limarí'
$ws.Range("S56").Value = 938
$ws.Range("T56").Value = 16
